$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend row 1 (header) with "sub.google.com" across columns P:BS (same style as O1 cascades
# automatically from the row's default style).
$ws.Range("P1:BS1").Value = "sub.google.com"

# Extend row 2 (data) with "x" across columns P:BS (same style as the other data cells in
# that row cascades automatically from the column's default style).
$ws.Range("P2:BS2").Value = "x"

# Reproduce the user's final selection: dragging/filling P1:BS2 leaves the active cell at the
# top-left of the fill (P1) with the full fill range selected.
$ws.Range("P1:BS2").Select()
